# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: force text storage for numeric-looking Price values
# by temporarily applying a text number format, then resetting the style
# back to Normal so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.715.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4662"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07915"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9767"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.743"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.965"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06942"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001005"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.628.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.336"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.127"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.076.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.790"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09366"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9425"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.330"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.347"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.352"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05850"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.77%  "

$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.869"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5654"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07360"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.220"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.25%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5335"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.139"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.852"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.354"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.035"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.82%  "
